$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.11017566666667
$ws.Range("H2").Value = 39.330527
$ws.Range("I2").Value = 0.1657114824704502
$ws.Range("J2").Value = 0.1657114824704501
$ws.Range("M2").Value = 19.98610666666667
$ws.Range("N2").Value = 59.95832
$ws.Range("O2").Value = 0.2969043109767812
$ws.Range("P2").Value = 0.2969043109767812
$ws.Range("Q2").Value = 262.0213692927379
$ws.Range("R2").Value = 2358.19232363464
$ws.Range("S2").Value = 0.04920045352382996
$ws.Range("T2").Value = 0.04920045352382996
$ws.Range("G3").Value = 13.11017566666667
$ws.Range("H3").Value = 39.330527
$ws.Range("I3").Value = 0.1657114824704502
$ws.Range("J3").Value = 0.1657114824704501
$ws.Range("O3").Value = 0.4664722083712238
$ws.Range("P3").Value = 0.4664722083712239
$ws.Range("Q3").Value = 411.6669319227021
$ws.Range("R3").Value = 3705.002387304319
$ws.Range("S3").Value = 0.07729980118046023
$ws.Range("T3").Value = 0.07729980118046023
$ws.Range("G4").Value = 13.11017566666667
$ws.Range("H4").Value = 39.330527
$ws.Range("I4").Value = 0.1657114824704502
$ws.Range("J4").Value = 0.1657114824704501
$ws.Range("O4").Value = 0.236623480651995
$ws.Range("P4").Value = 0.236623480651995
$ws.Range("Q4").Value = 208.8228635120693
$ws.Range("R4").Value = 1879.405771608624
$ws.Range("S4").Value = 0.03921122776615997
$ws.Range("T4").Value = 0.03921122776615997
$ws.Range("G5").Value = 51.42568199999999
$ws.Range("I5").Value = 0.6500161567583834
$ws.Range("J5").Value = 0.6500161567583833
$ws.Range("M5").Value = 19.98610666666667
$ws.Range("N5").Value = 59.95832
$ws.Range("O5").Value = 0.2969043109767812
$ws.Range("P5").Value = 0.2969043109767812
$ws.Range("Q5").Value = 1027.79916585808
$ws.Range("R5").Value = 9250.192492722719
$ws.Range("S5").Value = 0.1929925991461232
$ws.Range("T5").Value = 0.1929925991461232
$ws.Range("G6").Value = 51.42568199999999
$ws.Range("I6").Value = 0.6500161567583834
$ws.Range("J6").Value = 0.6500161567583833
$ws.Range("O6").Value = 0.4664722083712238
$ws.Range("P6").Value = 0.4664722083712239
$ws.Range("S6").Value = 0.3032144721200587
$ws.Range("T6").Value = 0.3032144721200587
$ws.Range("G7").Value = 51.42568199999999
$ws.Range("I7").Value = 0.6500161567583834
$ws.Range("J7").Value = 0.6500161567583833
$ws.Range("O7").Value = 0.236623480651995
$ws.Range("P7").Value = 0.236623480651995
$ws.Range("Q7").Value = 819.1238963033279
$ws.Range("R7").Value = 7372.115066729951
$ws.Range("S7").Value = 0.1538090854922015
$ws.Range("T7").Value = 0.1538090854922015
$ws.Range("I8").Value = 0.1842723607711665
$ws.Range("J8").Value = 0.1842723607711665
$ws.Range("M8").Value = 19.98610666666667
$ws.Range("N8").Value = 59.95832
$ws.Range("O8").Value = 0.2969043109767812
$ws.Range("P8").Value = 0.2969043109767812
$ws.Range("Q8").Value = 291.3696478496977
$ws.Range("R8").Value = 2622.32683064728
$ws.Range("S8").Value = 0.05471125830682804
$ws.Range("T8").Value = 0.05471125830682803
$ws.Range("I9").Value = 0.1842723607711665
$ws.Range("J9").Value = 0.1842723607711665
$ws.Range("O9").Value = 0.4664722083712238
$ws.Range("P9").Value = 0.4664722083712239
$ws.Range("S9").Value = 0.08595793507070491
$ws.Range("T9").Value = 0.08595793507070491
$ws.Range("I10").Value = 0.1842723607711665
$ws.Range("J10").Value = 0.1842723607711665
$ws.Range("O10").Value = 0.236623480651995
$ws.Range("P10").Value = 0.236623480651995
$ws.Range("S10").Value = 0.04360316739363355
$ws.Range("T10").Value = 0.04360316739363355
